# Daily attendance processing - 2025-10-23 19:42:26
# Normalize the "Recorded By" (column G) entries: for any cell whose value is a
# comma-separated list of two or more recorders, swap the last two entries
# (keeping any earlier entries in place). Single-value cells are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val -split ",\s*"

    if ($parts.Count -ge 2) {
        $last = $parts.Count - 1
        $tmp = $parts[$last]
        $parts[$last] = $parts[$last - 1]
        $parts[$last - 1] = $tmp
        $cell.Value = [string]::Join(", ", $parts)
    }
}
